# Apply "Added vocabulary concept ids" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Expand / rename a few vocabulary_name values to their fuller descriptions
$ws.Range("C48").Value = "Clinical Classifications Software for ICD-9-CM (HCUP)"
$ws.Range("C50").Value = "Gemscript NHS dictionary of medicine and devices (NHS)"
$ws.Range("C51").Value = "Hospital Episode Statistics Specialty (NHS)"
$ws.Range("C59").Value = "International Currency Symbol (ISO 4217)"

# 2) Add a new column G with concept ids, one per data row (rows 2-59)
$ws.Range("G1").Value = "CONCEPT_ID"

$conceptId = 44819096
for ($r = 2; $r -le 59; $r++) {
    $ws.Cells.Item($r, 7).Value = $conceptId
    $conceptId = $conceptId + 1
}

# 3) Update the view state to match the saved selection / scroll position
$ws.Range("I23").Select()
$excel.ActiveWindow.ScrollRow = 20
